$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 96: TRADING_ATTEMPT for TRX
$r = 96
$ws.Cells.Item($r, 1).Value = "2025-10-23T12:43:39.553291"
$ws.Cells.Item($r, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item($r, 3).Value = "TRX"
$ws.Cells.Item($r, 4).Value = "UNKNOWN"
$ws.Cells.Item($r, 5).Value = 0.3220488363629825
$ws.Cells.Item($r, 6).Font.Bold = $false
$ws.Cells.Item($r, 7).Font.Bold = $false
$ws.Cells.Item($r, 8).Font.Bold = $false
$ws.Cells.Item($r, 9).Font.Bold = $false
$ws.Cells.Item($r, 10).Font.Bold = $false
$ws.Cells.Item($r, 11).Value = "ATTEMPT"
$ws.Cells.Item($r, 12).Value = "Attempting trade 1/1"

# Row 97: POSITION_OPENED for TRX
$r = 97
$ws.Cells.Item($r, 1).Value = "2025-10-23T12:43:41.235374"
$ws.Cells.Item($r, 2).Value = "POSITION_OPENED"
$ws.Cells.Item($r, 3).Value = "TRX"
$ws.Cells.Item($r, 4).Value = "UNKNOWN"
$ws.Cells.Item($r, 5).Value = 0.3220488363629825
$ws.Cells.Item($r, 6).Value = 90
$ws.Cells.Item($r, 7).Value = 1
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Font.Bold = $false
$ws.Cells.Item($r, 10).Font.Bold = $false
$ws.Cells.Item($r, 11).Value = "SUCCESS"
$ws.Cells.Item($r, 12).Font.Bold = $false
